$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix G3: correct the truncated adhesive type text (missing closing paren)
$ws.Range("G3").Value = "PSA - Adhesive Transfer Tape (ATT)"

# Add new row 4 of data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "oca2_p3"
$ws.Range("C4").Value = "oca2_dma.xml"
$ws.Range("D4").Value = "oca2_compression.xml"
$ws.Range("E4").Value = "oca2_tension.xlsx"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "PSA - Adhesive Transfer Tape (ATT)"
$ws.Range("H4").Value = "Phase 3 (DMA + Compression + Tension)"
$ws.Range("I4").Value = "Tension"
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = "oca2_p3_ModelFiles.zip"

# Column width adjustments (best fit widths for columns D and H)
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668
$ws.Columns.Item(8).ColumnWidth = 33.333333333333336

# Update selection to match target state
$ws.Range("J16").Select()
